$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update mat_rank (matrices score, column F) values for all data rows
$ws.Range("F2").Value = 13.25581603006527
$ws.Range("F3").Value = 13.03564410204013
$ws.Range("F4").Value = 8.387878449008936
$ws.Range("F5").Value = 8.356292063322577
$ws.Range("F6").Value = 8.180015286402934
$ws.Range("F7").Value = 5.393336665672788

# Row 8 (A8=6) and Row 9 (A9=7) swap places in the ranking:
# the record that was in row 9 (prolificid 30 / 60d5775a99b502eec8cf56b4 / Shadaisia)
# moves up to row 8, and vice versa.
$ws.Range("B8").Value = 30
$ws.Range("C8").Value = "60d5775a99b502eec8cf56b4"
$ws.Range("D8").Value = "Shadaisia"
$ws.Range("F8").Value = 5.339669197139461
$ws.Range("G8").Value = "Black or African American"

$ws.Range("B9").Value = 32
$ws.Range("C9").Value = "6036f9b3b1842f8b659b18c7"
$ws.Range("D9").Value = "Kellie"
$ws.Range("F9").Value = 5.108019693417147
$ws.Range("G9").Value = "White"

$ws.Range("F10").Value = 4.268640122598316
$ws.Range("F11").Value = 4.254495598246366
$ws.Range("F12").Value = 2.420025270519735
$ws.Range("F13").Value = 1.496024677253027
